# Update Name of Algo
# Applies the updated KNN imputation results to column B and D of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = -7.271000000000001
    "B3"  = 5.617
    "B14" = 5.755999999999999
    "B16" = 5.375
    "D18" = -8.225000000000001
    "B21" = 8.996
    "B23" = 7.598000000000001
    "D24" = -7.555
    "B25" = 6.334000000000001
    "D25" = -7.792
    "B26" = 6.637
    "D27" = -8.279
    "B29" = 5.689
    "D30" = -7.305
    "D31" = -7.930000000000001
    "D39" = -7.782999999999999
    "B40" = 9.137
    "D42" = -8.193000000000001
    "D48" = -7.007
    "D51" = -8.272
    "D52" = -8.1
    "B53" = 5.805
    "D55" = -8.225
    "D56" = -8.306000000000001
    "B57" = 4.981999999999999
    "D57" = -8.17
    "B59" = 5.145000000000001
    "D60" = -8.279
    "B65" = 5.997000000000001
    "B69" = 5.164
    "D73" = -8.326000000000002
    "D74" = -7.888
    "B79" = 6.406000000000001
    "B83" = 5.553000000000001
    "D89" = -6.351
    "D90" = -7.472
    "B91" = 5.469
    "D92" = -6.35
    "B93" = 5.635
    "B100" = 5.396
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
